$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the values in column B for rows 66 to 82 from "n" to "y"
$ws.Range("B66:B82").Value = "y"

# Update the sheet view: set topLeftCell and change the selection
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B66:B82").Select()
$excel.ActiveCell = $ws.Range("B66")
